$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.231.36"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.784.23"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'224.62"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'31.86"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "2.039.40"
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").Value = "'11.16"
$ws.Range("E13").Value = "  +6.95%  "
$ws.Range("D14").Value = "1.777.58"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").Value = "34.228.44"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'4.21"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'68.74"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "'254.62"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'10.36"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'4.20"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "'157.33"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'16.40"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "'7.01"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("D35").Value = "1.440.21"
$ws.Range("E35").Value = "  -7.33%  "
$ws.Range("E36").Value = "  -3.33%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "'0.890"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("D44").Value = "'0.0510"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'5.83"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.939.81"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").Value = "'12.18"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'98.57"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "'49.46"
$ws.Range("E51").Value = "  -7.06%  "
